$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values (columns H through AA) per the updated election results
$ws.Range("H2").Value  = 25
$ws.Range("I2").Value  = 62
$ws.Range("J2").Value  = 339
$ws.Range("K2").Value  = 3
$ws.Range("L2").Value  = 101
$ws.Range("M2").Value  = 2
$ws.Range("N2").Value  = 64
$ws.Range("O2").Value  = 0
$ws.Range("P2").Value  = 3
$ws.Range("Q2").Value  = 1
$ws.Range("R2").Value  = 3
$ws.Range("S2").Value  = 42
$ws.Range("T2").Value  = 44
$ws.Range("U2").Value  = 3
$ws.Range("V2").Value  = 507
$ws.Range("W2").Value  = 1
$ws.Range("X2").Value  = 509
$ws.Range("Y2").Value  = 2
$ws.Range("Z2").Value  = 5
$ws.Range("AA2").Value = 1

$wb.Save()
